# Reapply "Merge remote-tracking branch 'origin/Dev_0.0.1' into ArtWork"
# - removes the old row 8 (leftover monster entry)
# - fixes the data-type marker in E3 (long -> float)
# - renames the monster placeholder names to mon1..mon4
# - rebalances the spawn distance / count numbers
# - renames the movement-pattern values (left/right -> Left/Right/LeftFast)
# - fills in the new per-row Key/Score/Lood/MoveKey data in columns G,H,K,L
# - adds a brand-new "Prefabs" column (M) with its header/description/type/data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Drop the old trailing row (monster #5 leftover row) - shifts nothing else
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Delete()

# ---------------------------------------------------------------------------
# Row 3 (data-type marker row): E3 switches from "long" to "float"
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "float"

# ---------------------------------------------------------------------------
# Row 4 (mon1 / 잼민이 -> mon1)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "mon1"
$ws.Range("C4").Value = 5000
$ws.Range("D4").Value = 400
$ws.Range("G4").Value = "mon1"
$ws.Range("L4").Value = "Left"

# ---------------------------------------------------------------------------
# Row 5 (mon2 / 급식충 -> mon2)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "mon2"
$ws.Range("B5").Value = 10000
$ws.Range("C5").Value = 20000
$ws.Range("G5").Value = "mon2"
$ws.Range("H5").Value = 20
$ws.Range("K5").Value = "0,0,0"
$ws.Range("L5").Value = "Left"

# ---------------------------------------------------------------------------
# Row 6 (mon3 / 학식충 -> mon3)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "mon3"
$ws.Range("B6").Value = 50000
$ws.Range("C6").Value = 80000
$ws.Range("G6").Value = "mon3"
$ws.Range("H6").Value = 30
$ws.Range("K6").Value = "0,0,0"
$ws.Range("L6").Value = "Right"

# ---------------------------------------------------------------------------
# Row 7 (mon4 / 금태양 -> mon4)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "mon4"
$ws.Range("B7").Value = 15000
$ws.Range("C7").Value = 60000
$ws.Range("G7").Value = "mon4"
$ws.Range("H7").Value = 40
$ws.Range("K7").Value = "0,0,0"
$ws.Range("L7").Value = "LeftFast"

# ---------------------------------------------------------------------------
# New column M ("Prefabs") - copy the formatting from the existing column K
# header/description/type rows (s=2 / s=1 pattern), except M2 which uses the
# plain s=1 style instead of the s=3/s=4 used by the rest of row 2.
# ---------------------------------------------------------------------------
$ws.Range("K1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("K3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

$ws.Range("M1").Value = "Prefabs"
$ws.Range("M2").Value = "인게임 캐릭터 프리팹"
$ws.Range("M3").Value = "string"
$ws.Range("M4").Value = "mon1"
$ws.Range("M5").Value = "mon1"
$ws.Range("M6").Value = "mon1"
$ws.Range("M7").Value = "mon1"

$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# ---------------------------------------------------------------------------
# Selection moves to D5
# ---------------------------------------------------------------------------
$ws.Range("D5").Select()
